$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1 title text ---------------------------------------------------
$ws.Range("A1").Value = "fimra01"

# --- Replicate the row3/row4 formatting pattern down onto the new ----
# --- product blocks (rows 6-7, 8-9, 10-11) before filling values -----
$ws.Range("A3:E4").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)
$ws.Range("A3:E4").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)
$ws.Range("A3:E4").Copy()
$ws.Range("A10:E11").PasteSpecial(-4122)

# --- Merges for the two-row product blocks - do this BEFORE the final
# --- formatting pass below, because Merge() redraws the inner border
# --- (drops the shared edge) which the source file does not do.
$ws.Range("B6:B7").Merge()
$ws.Range("A6:A7").Merge()
$ws.Range("B8:B9").Merge()
$ws.Range("A8:A9").Merge()
$ws.Range("B10:B11").Merge()
$ws.Range("A10:A11").Merge()

# --- Re-stamp the original (un-redrawn) borders/format over the -------
# --- now-merged blocks, undoing Merge()'s automatic border surgery. ---
$ws.Range("A3:E4").Copy()
$ws.Range("A6:E7").PasteSpecial(-4122)
$ws.Range("A3:E4").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)
$ws.Range("A3:E4").Copy()
$ws.Range("A10:E11").PasteSpecial(-4122)

# --- Row 3 (was D2/27/12/czarny -> M1/ /21/styropian) -----------------
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = 21
$ws.Range("D3").Value = "styropian"

# --- Row 4 (C4 3 -> 32) ------------------------------------------------
$ws.Range("C4").Value = 32

# --- Row 6 (was M1/50/50/styropian -> M3/ /12/styropian) --------------
$ws.Range("A6").Value = "M3"
$ws.Range("B6").Value = $null
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = "styropian"

# --- Row 7 (new) --------------------------------------------------------
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = "czarny"

# --- Row 8 (new) --------------------------------------------------------
$ws.Range("A8").Value = "woodenStands"
$ws.Range("C8").Value = 33
$ws.Range("D8").Value = "biały"

# --- Row 9 (new) --------------------------------------------------------
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = "czarny"

# --- Row 10 (new) --------------------------------------------------------
$ws.Range("A10").Value = "stands"
$ws.Range("C10").Value = 12

# D10 must stay literal text "45" (not get coerced to a number) while
# keeping the D-column style (s=4). Build the text on a scratch cell far
# outside the used range, then bring just the *value* across - a values
# -only paste keeps the destination's existing style untouched.
$ws.Range("ZZ1").Value = "'45"
$ws.Range("ZZ1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# --- Row 11 (new) --------------------------------------------------------
$ws.Range("C11").Value = 3

$ws.Range("ZZ1").Value = "'90"
$ws.Range("ZZ1").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
